$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 1327.75
$ws.Cells.Item(138, 9).Value = 758
$ws.Cells.Item(138, 10).Value = 1821.5333
$ws.Cells.Item(138, 11).Value = 2274
$ws.Cells.Item(138, 12).Value = 5464.5999
$ws.Cells.Item(138, 13).Value = 2866
$ws.Cells.Item(138, 14).Value = -15744.5999
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(118, 8).Value = 18000
$ws.Cells.Item(118, 10).Value = 18000
$ws.Cells.Item(118, 12).Value = 18000
$ws.Cells.Item(118, 14).Value = -21314
$ws.Cells.Item(132, 8).Value = 3001.1538
$ws.Cells.Item(132, 9).Value = 2565.0908
$ws.Cells.Item(132, 10).Value = 5399.5
$ws.Cells.Item(132, 11).Value = 7695.2724
$ws.Cells.Item(132, 12).Value = 16198.5
$ws.Cells.Item(132, 13).Value = -5165.2724
$ws.Cells.Item(132, 14).Value = -21258.5
$ws.Cells.Item(134, 8).Value = 16129953
$ws.Cells.Item(134, 9).Value = 979
$ws.Cells.Item(134, 10).Value = 100000620
$ws.Cells.Item(134, 11).Value = 2937
$ws.Cells.Item(134, 12).Value = 300001860
$ws.Cells.Item(134, 13).Value = -402
$ws.Cells.Item(134, 14).Value = -300006930
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 20900.8
$ws.Cells.Item(5, 9).Value = 999.5
$ws.Cells.Item(5, 11).Value = 2998.5
$ws.Cells.Item(5, 13).Value = -2886.5
$ws.Cells.Item(7, 8).Value = 308
$ws.Cells.Item(7, 9).Value = 0
$ws.Cells.Item(7, 10).Value = 308
$ws.Cells.Item(7, 11).Value = 0
$ws.Cells.Item(7, 12).Value = 924
$ws.Cells.Item(7, 13).ClearContents()
$ws.Cells.Item(7, 14).Value = -1148
$ws.Cells.Item(12, 8).Value = 73.76470999999999
$ws.Cells.Item(12, 9).Value = 204.375
$ws.Cells.Item(12, 10).Value = 33.576923
$ws.Cells.Item(12, 11).Value = 613.125
$ws.Cells.Item(12, 12).Value = 100.730769
$ws.Cells.Item(12, 13).Value = -440.125
$ws.Cells.Item(12, 14).Value = -446.730769
$ws.Cells.Item(20, 8).Value = 970
$ws.Cells.Item(20, 9).Value = 197.5
$ws.Cells.Item(20, 10).Value = 2000
$ws.Cells.Item(20, 11).Value = 592.5
$ws.Cells.Item(20, 12).Value = 6000
$ws.Cells.Item(20, 13).Value = -365.5
$ws.Cells.Item(20, 14).Value = -6454
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(21, 9).Value = 0
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).ClearContents()
$ws.Cells.Item(21, 14).ClearContents()
$ws.Cells.Item(22, 8).Value = 1579.8
$ws.Cells.Item(22, 9).Value = 999.5
$ws.Cells.Item(22, 10).Value = 1966.6666
$ws.Cells.Item(22, 11).Value = 2998.5
$ws.Cells.Item(22, 12).Value = 5899.9998
$ws.Cells.Item(22, 13).Value = -2829.5
$ws.Cells.Item(22, 14).Value = -6237.9998
$ws.Cells.Item(27, 8).Value = 1579.8
$ws.Cells.Item(27, 9).Value = 999.5
$ws.Cells.Item(27, 10).Value = 1966.6666
$ws.Cells.Item(27, 11).Value = 2998.5
$ws.Cells.Item(27, 12).Value = 5899.9998
$ws.Cells.Item(27, 13).Value = -2896.5
$ws.Cells.Item(27, 14).Value = -6103.9998
$ws.Cells.Item(34, 8).Value = 49383264
$ws.Cells.Item(34, 9).Value = 500000030
$ws.Cells.Item(34, 10).Value = 13333924
$ws.Cells.Item(34, 11).Value = 1500000090
$ws.Cells.Item(34, 12).Value = 40001772
$ws.Cells.Item(34, 13).Value = -1500000006
$ws.Cells.Item(34, 14).Value = -40001940
$ws.Cells.Item(40, 8).Value = 285
$ws.Cells.Item(40, 9).Value = 285
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 1140
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -1071
$ws.Cells.Item(40, 14).ClearContents()
$ws.Cells.Item(46, 8).Value = 1727.4
$ws.Cells.Item(46, 9).Value = 350
$ws.Cells.Item(46, 10).Value = 1939.3077
$ws.Cells.Item(46, 11).Value = 1050
$ws.Cells.Item(46, 12).Value = 5817.9231
$ws.Cells.Item(46, 13).Value = -959
$ws.Cells.Item(46, 14).Value = -5999.9231
$ws.Cells.Item(64, 8).Value = 2024249.9
$ws.Cells.Item(64, 10).Value = 2024249.9
$ws.Cells.Item(64, 12).Value = 6072749.699999999
$ws.Cells.Item(64, 14).Value = -6073289.699999999
$ws.Cells.Item(67, 8).Value = 2024249.9
$ws.Cells.Item(67, 10).Value = 2024249.9
$ws.Cells.Item(67, 12).Value = 6072749.699999999
$ws.Cells.Item(67, 14).Value = -6074621.699999999
$ws.Cells.Item(75, 8).Value = 3844.5
$ws.Cells.Item(75, 10).Value = 4576.364
$ws.Cells.Item(75, 12).Value = 13729.092
$ws.Cells.Item(75, 14).Value = -15725.092
$ws.Cells.Item(78, 8).Value = 3844.5
$ws.Cells.Item(78, 10).Value = 4576.364
$ws.Cells.Item(78, 12).Value = 41187.276
$ws.Cells.Item(78, 14).Value = -51171.276
$ws.Cells.Item(94, 8).Value = 2120
$ws.Cells.Item(94, 9).Value = 490
$ws.Cells.Item(94, 10).Value = 3750
$ws.Cells.Item(94, 11).Value = 1470
$ws.Cells.Item(94, 12).Value = 11250
$ws.Cells.Item(94, 13).Value = -794
$ws.Cells.Item(94, 14).Value = -12602
$ws.Cells.Item(100, 8).Value = 41340.332
$ws.Cells.Item(100, 10).Value = 41340.332
$ws.Cells.Item(100, 12).Value = 124020.996
$ws.Cells.Item(100, 14).Value = -125642.996
$ws.Cells.Item(108, 8).Value = 1879.091
$ws.Cells.Item(108, 9).Value = 524
$ws.Cells.Item(108, 10).Value = 3008.3333
$ws.Cells.Item(108, 11).Value = 1572
$ws.Cells.Item(108, 12).Value = 9024.999899999999
$ws.Cells.Item(108, 13).Value = 1308
$ws.Cells.Item(108, 14).Value = -14784.9999
$ws.Cells.Item(122, 8).Value = 348.56
$ws.Cells.Item(122, 9).Value = 184.38889
$ws.Cells.Item(122, 10).Value = 770.7143
$ws.Cells.Item(122, 11).Value = 1659.50001
$ws.Cells.Item(122, 12).Value = 6936.428699999999
$ws.Cells.Item(122, 13).Value = 790.49999
$ws.Cells.Item(122, 14).Value = -11836.4287
$ws.Cells.Item(135, 8).Value = 20900.8
$ws.Cells.Item(135, 9).Value = 999.5
$ws.Cells.Item(135, 11).Value = 8995.5
$ws.Cells.Item(135, 13).Value = -6460.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(68, 8).Value = 1842.3125
$ws.Cells.Item(68, 9).Value = 1680.7273
$ws.Cells.Item(68, 10).Value = 2197.8
$ws.Cells.Item(68, 11).Value = 1680.7273
$ws.Cells.Item(68, 12).Value = 2197.8
$ws.Cells.Item(68, 13).Value = -931.7273
$ws.Cells.Item(68, 14).Value = -3695.8
$ws.Cells.Item(71, 8).Value = 1842.3125
$ws.Cells.Item(71, 9).Value = 1680.7273
$ws.Cells.Item(71, 10).Value = 2197.8
$ws.Cells.Item(71, 11).Value = 8403.636500000001
$ws.Cells.Item(71, 12).Value = 10989
$ws.Cells.Item(71, 13).Value = -4659.636500000001
$ws.Cells.Item(71, 14).Value = -18477
$ws.Cells.Item(136, 8).Value = 557149.8
$ws.Cells.Item(136, 9).Value = 626043.5600000001
$ws.Cells.Item(136, 11).Value = 1878130.68
$ws.Cells.Item(136, 13).Value = -1875580.68
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 5961.1924
$ws.Cells.Item(132, 9).Value = 1687.375
$ws.Cells.Item(132, 10).Value = 12799.3
$ws.Cells.Item(132, 11).Value = 5062.125
$ws.Cells.Item(132, 12).Value = 38397.89999999999
$ws.Cells.Item(132, 13).Value = -2532.125
$ws.Cells.Item(132, 14).Value = -43457.89999999999
